$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.652167666666667
$ws.Range("H2").Value = 16.956503
$ws.Range("I2").Value = 0.1860329065948871
$ws.Range("J2").Value = 0.1860329065948871
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.672264666666667
$ws.Range("N2").Value = 8.016794000000001
$ws.Range("O2").Value = 0.06772620019093417
$ws.Range("P2").Value = 0.06772620019093417
$ws.Range("Q2").Value = 15.10408794570911
$ws.Range("R2").Value = 135.936791511382
$ws.Range("S2").Value = 0.01259930187414668
$ws.Range("T2").Value = 0.01259930187414668
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.652167666666667
$ws.Range("H3").Value = 16.956503
$ws.Range("I3").Value = 0.1860329065948871
$ws.Range("J3").Value = 0.1860329065948871
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.41886
$ws.Range("N3").Value = 82.25658
$ws.Range("O3").Value = 0.6949069171668364
$ws.Range("P3").Value = 0.6949069171668364
$ws.Range("Q3").Value = 154.97599394886
$ws.Range("R3").Value = 1394.78394553974
$ws.Range("S3").Value = 0.129275553613439
$ws.Range("T3").Value = 0.129275553613439
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.652167666666667
$ws.Range("H4").Value = 16.956503
$ws.Range("I4").Value = 0.1860329065948871
$ws.Range("J4").Value = 0.1860329065948871
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.365757
$ws.Range("N4").Value = 28.097271
$ws.Range("O4").Value = 0.2373668826422294
$ws.Range("P4").Value = 0.2373668826422294
$ws.Range("Q4").Value = 52.93682888925701
$ws.Range("R4").Value = 476.431460003313
$ws.Range("S4").Value = 0.04415805110730138
$ws.Range("T4").Value = 0.04415805110730138
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.407289
$ws.Range("H5").Value = 49.221867
$ws.Range("I5").Value = 0.5400221369958743
$ws.Range("J5").Value = 0.5400221369958743
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.672264666666667
$ws.Range("N5").Value = 8.016794000000001
$ws.Range("O5").Value = 0.06772620019093417
$ws.Range("P5").Value = 0.06772620019093417
$ws.Range("Q5").Value = 43.84461867048868
$ws.Range("R5").Value = 394.601568034398
$ws.Range("S5").Value = 0.03657364735771866
$ws.Range("T5").Value = 0.03657364735771866
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.407289
$ws.Range("H6").Value = 49.221867
$ws.Range("I6").Value = 0.5400221369958743
$ws.Range("J6").Value = 0.5400221369958743
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.41886
$ws.Range("N6").Value = 82.25658
$ws.Range("O6").Value = 0.6949069171668364
$ws.Range("P6").Value = 0.6949069171668364
$ws.Range("Q6").Value = 449.86916007054
$ws.Range("R6").Value = 4048.82244063486
$ws.Range("S6").Value = 0.37526511842165
$ws.Range("T6").Value = 0.37526511842165
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.407289
$ws.Range("H7").Value = 49.221867
$ws.Range("I7").Value = 0.5400221369958743
$ws.Range("J7").Value = 0.5400221369958743
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.365757
$ws.Range("N7").Value = 28.097271
$ws.Range("O7").Value = 0.2373668826422294
$ws.Range("P7").Value = 0.2373668826422294
$ws.Range("Q7").Value = 153.666681802773
$ws.Range("R7").Value = 1383.000136224957
$ws.Range("S7").Value = 0.1281833712165056
$ws.Range("T7").Value = 0.1281833712165056
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.323166333333333
$ws.Range("H8").Value = 24.969499
$ws.Range("I8").Value = 0.2739449564092387
$ws.Range("J8").Value = 0.2739449564092387
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.672264666666667
$ws.Range("N8").Value = 8.016794000000001
$ws.Range("O8").Value = 0.06772620019093417
$ws.Range("P8").Value = 0.06772620019093417
$ws.Range("Q8").Value = 22.24170330735623
$ws.Range("R8").Value = 200.175329766206
$ws.Range("S8").Value = 0.01855325095906883
$ws.Range("T8").Value = 0.01855325095906884
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.323166333333333
$ws.Range("H9").Value = 24.969499
$ws.Range("I9").Value = 0.2739449564092387
$ws.Range("J9").Value = 0.2739449564092387
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.41886
$ws.Range("N9").Value = 82.25658
$ws.Range("O9").Value = 0.6949069171668364
$ws.Range("P9").Value = 0.6949069171668364
$ws.Range("Q9").Value = 228.21173245038
$ws.Range("R9").Value = 2053.90559205342
$ws.Range("S9").Value = 0.1903662451317474
$ws.Range("T9").Value = 0.1903662451317475
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.323166333333333
$ws.Range("H10").Value = 24.969499
$ws.Range("I10").Value = 0.2739449564092387
$ws.Range("J10").Value = 0.2739449564092387
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.365757
$ws.Range("N10").Value = 28.097271
$ws.Range("O10").Value = 0.2373668826422294
$ws.Range("P10").Value = 0.2373668826422294
$ws.Range("Q10").Value = 77.952753348581
$ws.Range("R10").Value = 701.5747801372289
$ws.Range("S10").Value = 0.06502546031842241
$ws.Range("T10").Value = 0.06502546031842242
